$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 425, shifting existing rows 425.. down by one.
$ws.Rows.Item(425).Insert()

# Populate the newly inserted row 425 with the new weekly price record.
$ws.Range("A425").Value = 5
$ws.Range("B425").Value = "Macroferia Regional de Talca"
$ws.Range("C425").Value = "Maule"
$ws.Range("D425").Value = 44931
$ws.Range("D425").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("E425").Value = 7
$ws.Range("F425").Value = 100114014
$ws.Range("G425").Value = "Betarraga"
$ws.Range("H425").Value = "Sin especificar"
$ws.Range("I425").Value = "Primera"
$ws.Range("J425").Value = 5000
$ws.Range("K425").Value = 700
$ws.Range("L425").Value = 700
$ws.Range("M425").Value = 700
$ws.Range("N425").Value = "`$/paquete 5 unidades"
$ws.Range("O425").Value = "Región del Maule"
$ws.Range("P425").Value = 140
$ws.Range("Q425").Value = 5
$ws.Range("R425").Value = "Hortaliza"
